$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A2").Value = "Pencil"
$ws2.Range("A3").Value = "Iphone"
$ws2.Range("A4").Value = "Toys"
$ws2.Range("A5").Value = "women's clothing"
$ws2.Range("A6").Value = "shoes"
$ws2.Range("A1").Value = "searchItem"
